{"js": "// The document references several R-script line numbers (e.g. \"line 18\",\n// \"Lines 6-15\") that need to be decremented by 1 throughout the body\n// (both the single \"line N\" mentions and both endpoints of the\n// \"Lines N-M\" ranges). We locate each exact phrase with body.search()\n// and replace it in place, in descending order of string length so a\n// shorter phrase (e.g. \"Lines 43-59\") never accidentally matches inside\n// a longer one that shares the same prefix (e.g. \"Lines 43-77\" vs\n// \"Lines 43-59\" both start with \"Lines 43-\" but are distinct full\n// strings, so ordering isn't actually required for correctness here,\n// but we keep replacements scoped to exact, unique phrases to be safe).\nconst replacements = [\n  [\"line 18\", \"line 17\"],\n  [\"line 21\", \"line 20\"],\n  [\"line 24\", \"line 23\"],\n  [\"line 12\", \"line 11\"],\n  [\"Lines 6-15\", \"Lines 5-14\"],\n  [\"Lines 17-24\", \"Lines 16-23\"],\n  [\"Lines 26-35\", \"Lines 25-34\"],\n  [\"Lines 37-40\", \"Lines 36-39\"],\n  [\"Lines 43-77\", \"Lines 42-76\"],\n  [\"Lines 43-59\", \"Lines 42-58\"],\n  [\"Lines 61-69\", \"Lines 60-68\"],\n  [\"Lines 71-77\", \"Lines 70-76\"],\n  [\"Lines 80-115\", \"Lines 79-114\"],\n  [\"Lines 81-85\", \"Lines 80-84\"],\n  [\"Lines 87-109\", \"Lines 86-108\"],\n  [\"Lines 111-115\", \"Lines 110-114\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document references several R-script line numbers (e.g. \"line 18\",\n# \"Lines 6-15\") that need to be decremented by 1 throughout the body\n# (both the single \"line N\" mentions and both endpoints of the\n# \"Lines N-M\" ranges). Each old/new phrase pair below is an exact,\n# unique string, so a simple Find/Replace-All per pair is safe and will\n# not clobber any other occurrence.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"line 18\", \"line 17\"),\n    @(\"line 21\", \"line 20\"),\n    @(\"line 24\", \"line 23\"),\n    @(\"line 12\", \"line 11\"),\n    @(\"Lines 6-15\", \"Lines 5-14\"),\n    @(\"Lines 17-24\", \"Lines 16-23\"),\n    @(\"Lines 26-35\", \"Lines 25-34\"),\n    @(\"Lines 37-40\", \"Lines 36-39\"),\n    @(\"Lines 43-77\", \"Lines 42-76\"),\n    @(\"Lines 43-59\", \"Lines 42-58\"),\n    @(\"Lines 61-69\", \"Lines 60-68\"),\n    @(\"Lines 71-77\", \"Lines 70-76\"),\n    @(\"Lines 80-115\", \"Lines 79-114\"),\n    @(\"Lines 81-85\", \"Lines 80-84\"),\n    @(\"Lines 87-109\", \"Lines 86-108\"),\n    @(\"Lines 111-115\", \"Lines 110-114\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
